$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - First Board Design
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 0.9

# Row 9 - Final Proposal
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 7

# Row 10 - Order Parts
$ws.Range("C10").Value = 10
$ws.Range("E10").Value = 10

# Row 11 - Graphing Software
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = 22
$ws.Range("E11").Value = 9
$ws.Range("G11").Value = 0.05

# Row 12 - Final Proposal Presentation
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 0.5

# Move the active selection to H2 (was C8)
$ws.Range("H2").Select()

# Force recalculation so the volatile WEEKNUM/TODAY formula in H2 refreshes
$excel.Calculate()
